$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = 43775
$ws.Range("B13").Value = "Requirements + RASD"
$ws.Range("C13").Value = 3
$ws.Range("C14").Formula = "=SUM(C4:C13)"
